$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 111, pushing existing rows 111-178 down to 112-179.
$ws.Rows(111).Insert()

# Populate the newly inserted row 111 with its data.
$ws.Cells.Item(111, 1).Value = 7
$ws.Cells.Item(111, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(111, 3).Value = "Ñuble"
$ws.Cells.Item(111, 4).Value = 44529
$ws.Cells.Item(111, 5).Value = 16
$ws.Cells.Item(111, 6).Value = 100112043
$ws.Cells.Item(111, 7).Value = "Pepino ensalada"
$ws.Cells.Item(111, 8).Value = "Sin especificar"
$ws.Cells.Item(111, 9).Value = "Primera"
$ws.Cells.Item(111, 10).Value = 100
$ws.Cells.Item(111, 11).Value = 8000
$ws.Cells.Item(111, 12).Value = 8500
$ws.Cells.Item(111, 13).Value = 8250
$ws.Cells.Item(111, 14).Value = "$/caja 80 unidades"
$ws.Cells.Item(111, 15).Value = "Región del Maule"
$ws.Cells.Item(111, 16).Value = 103
$ws.Cells.Item(111, 17).Value = 80
$ws.Cells.Item(111, 18).Value = "Hortaliza"

# Ensure the date cell keeps the date number format used by the other date cells.
$ws.Cells.Item(111, 4).NumberFormat = $ws.Cells.Item(112, 4).NumberFormat
